$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right before "总计", reproducing
#    the sheetId / relationship-id ordering Excel would naturally
#    assign (2022-Q1 -> sheetId 3 / rId3, 总计 -> sheetId 4 / rId4):
#    temporarily rename the existing 总计 sheet, add a fresh sheet
#    right after it (the fresh sheet claims the next sheetId, 4),
#    rename the fresh sheet to 总计, and rename the original
#    (now-repurposed) sheet to 2022-Q1. The repurposed sheet already
#    carries the old 总计 data/styles, which we overwrite below; the
#    fresh sheet starts empty and we populate it as the new 总计.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(3)
$total.Name = "TEMP_TOTAL_SHEET"
$freshTotal = $wb.Worksheets.Add($null, $total)
$freshTotal.Name = "总计"
$q1 = $total
$q1.Name = "2022-Q1"

# ------------------------------------------------------------------
# Helpers
# ------------------------------------------------------------------

# Write a string that Excel's grid would otherwise auto-convert to a
# number (fund codes with leading zeros, decimal-looking percentages,
# etc.) by pre-formatting the cell as Text, so it round-trips as the
# literal string instead of a numeric value.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Copy another cell's number format / style onto $cell (used to keep
# the index/row-number column matching the sheet's existing "s=2"
# style instead of picking up the freshly-minted text style above).
function Copy-Style($srcCell, $dstCell) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
}

# ------------------------------------------------------------------
# 2) "2022-Q1" sheet content
# ------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item(2).Cells.Item(2, 1)

Set-TextValue $q1.Cells.Item(1, 2) "基金代码"
Set-TextValue $q1.Cells.Item(1, 3) "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"
Copy-Style $styleSrc $q1.Cells.Item(1, 2)
Copy-Style $styleSrc $q1.Cells.Item(1, 3)

$q1Data = @(
    @("005457", "景顺长城量化小盘股票",         "9.49", "93.39", "1.91", "0.1813", 7),
    @("005313", "万家中证1000指数增强A",        "9.01", "93.72", "1.21", "0.1090", 4),
    @("005314", "万家中证1000指数增强C",        "4.95", "93.72", "1.21", "0.0599", 4),
    @("620007", "金元顺安优质精选灵活配置混合A", "0.75", "39.31", "2.08", "0.0156", 5),
    @("011179", "浙商智选食品饮料股票A",         "0.22", "91.35", "6.63", "0.0146", 5),
    @("001375", "金元顺安优质精选灵活配置混合C", "0.69", "39.31", "2.08", "0.0144", 5),
    @("011180", "浙商智选食品饮料股票C",         "0.05", "91.35", "6.63", "0.0033", 5)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    Copy-Style $styleSrc $q1.Cells.Item($r, 1)
    Set-TextValue $q1.Cells.Item($r, 2) $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    Set-TextValue $q1.Cells.Item($r, 4) $row[2]
    Set-TextValue $q1.Cells.Item($r, 5) $row[3]
    Set-TextValue $q1.Cells.Item($r, 6) $row[4]
    Set-TextValue $q1.Cells.Item($r, 7) $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 3) "总计" sheet content: previous 3 rows with a new 2022-Q1 row
#    inserted at the top.
# ------------------------------------------------------------------
$freshTotal.Cells.Item(1, 2).Value = "日期"
$freshTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$freshTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 7, 0.4),
    @("2021-Q4", 2, 0.01),
    @("2021-Q3", 2, 0.02)
)

$r = 2
foreach ($row in $totalData) {
    $freshTotal.Cells.Item($r, 1).Value = $r - 2
    Copy-Style $styleSrc $freshTotal.Cells.Item($r, 1)
    Set-TextValue $freshTotal.Cells.Item($r, 2) $row[0]
    $freshTotal.Cells.Item($r, 3).Value = $row[1]
    $freshTotal.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
